$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit is effectively a cyclic shift of the weekly price-record rows:
#   row2 <- old row3, row3 <- old row4, row4 <- old row2
# across columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg).
$cols = @("D", "M", "N", "O", "P", "S")

# Capture the original values first (use Value2 for a clean read).
$orig2 = @{}
$orig3 = @{}
$orig4 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range($col + "2").Value2
    $orig3[$col] = $ws.Range($col + "3").Value2
    $orig4[$col] = $ws.Range($col + "4").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "2").Value = $orig3[$col]
    $ws.Range($col + "3").Value = $orig4[$col]
    $ws.Range($col + "4").Value = $orig2[$col]
}
